$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: new test result values (verifying 4 values)
$ws.Range("A2").Value = "test_tour_cashweekly"

$ws.Range("B2").Formula = '="$658.7"'
$ws.Range("B2").Copy()
$ws.Range("B2").PasteSpecial(-4163)

$ws.Range("C2").Formula = '="$4,279.8"'
$ws.Range("C2").Copy()
$ws.Range("C2").PasteSpecial(-4163)

$ws.Range("D2").NumberFormat = "0.0%"
$ws.Range("D2").Value = 0.911

$ws.Range("E2").Formula = '="$53.2"'
$ws.Range("E2").Copy()
$ws.Range("E2").PasteSpecial(-4163)

$ws.Range("F2").ClearContents()
$ws.Range("G2").ClearContents()
$ws.Range("H2").ClearContents()
$ws.Range("I2").ClearContents()
$ws.Range("J2").ClearContents()
$ws.Range("K2").ClearContents()

# Rows 3-5: test names (unchanged text, but shared-string table got reshuffled)
$ws.Range("A3").Value = "test_tour_revenue"
$ws.Range("A4").Value = "test_tour_inventory"
$ws.Range("A5").Value = "test_target_audience_builder"

# Column A needs to widen to fit the longest new test name
$ws.Columns.Item(1).AutoFit()
$excel.CutCopyMode = $false

# Selection moved to D8
$ws.Range("D8").Select() | Out-Null
